# [Salesmanager] Edit Product - Fixbug
# Update a handful of price/quantity/sold figures on the product import
# sheet and remove the three trailing rows (42-44) that were added by
# mistake while testing "add product" during the original session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quantity ("Số lượng") fixes
$ws.Range("D2").Value = 189
$ws.Range("D8").Value = 200
$ws.Range("D19").Value = 200

# Price ("Giá") fix
$ws.Range("C16").Value = 148000

# Row 40 (product #39) - price/quantity/sold corrections
$ws.Range("C40").Value = 1000
$ws.Range("D40").Value = 1
$ws.Range("J40").Value = 2

# Row 41 (product #40) - price/quantity/sold corrections
$ws.Range("C41").Value = 70000
$ws.Range("D41").Value = 1
$ws.Range("J41").Value = 2

# Drop the stray test rows 42-44 (sản phẩm mới 1214 / thêm cái nữa lúc 1233 / sp mới lúc 1257)
$ws.Range("A42:K44").EntireRow.Delete()

# Leave the selection on the new last data row, matching the saved view
[void]$ws.Range("C41").Select()
